$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing cell contents (keeps formatting, e.g. the bold header style),
# and lets the shared-strings table be rebuilt from only what is actually used.
$ws.Cells.ClearContents()

# Header row (unchanged text, just re-written since contents were cleared)
$headers = @("Sending cluster","Ligand symbol","Receptor symbol","Target cluster","Ligand-expressing cells","Ligand detection rate","Ligand average expression value","Ligand total expression value","Ligand derived specificity of average expression value","Ligand derived specificity of total expression value","Receptor-expressing cells","Receptor detection rate","Receptor average expression value","Receptor total expression value","Receptor derived specificity of average expression value","Receptor derived specificity of total expression value","Edge average expression weight","Edge total expression weight","Edge average expression derived specificity","Edge total expression derived specificity")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $ws.Cells.Item(1, $c + 1).Value2 = $headers[$c]
}

# Updated data rows (new TPM-derived values; table now has 6 data rows instead of 8)
$data = New-Object 'object[,]' 6,20
$data[0,0] = "ECs"
$data[0,1] = "Rspo3"
$data[0,2] = "Lgr5"
$data[0,3] = "ECs"
$data[0,4] = 2
$data[0,5] = 0.6666666666666666
$data[0,6] = 0.1043256666666667
$data[0,7] = 0.312977
$data[0,8] = 0.02547563162231953
$data[0,9] = 0.02547563162231953
$data[0,10] = 3
$data[0,11] = 1
$data[0,12] = 0.4967476666666666
$data[0,13] = 1.490243
$data[0,14] = 0.08283342158193596
$data[0,15] = 0.08283342158193596
$data[0,16] = 0.05182353149011111
$data[0,17] = 0.466411783411
$data[0,18] = 0.002110233734237693
$data[0,19] = 0.002110233734237693
$data[1,0] = "ECs"
$data[1,1] = "Rspo3"
$data[1,2] = "Lgr5"
$data[1,3] = "FAPs"
$data[1,4] = 2
$data[1,5] = 0.6666666666666666
$data[1,6] = 0.1043256666666667
$data[1,7] = 0.312977
$data[1,8] = 0.02547563162231953
$data[1,9] = 0.02547563162231953
$data[1,10] = 3
$data[1,11] = 1
$data[1,12] = 1.968976333333333
$data[1,13] = 5.906929
$data[1,14] = 0.3283297691125296
$data[1,15] = 0.3283297691125296
$data[1,16] = 0.2054147686258889
$data[1,17] = 1.848732917633
$data[1,18] = 0.00836440824855203
$data[1,19] = 0.00836440824855203
$data[2,0] = "ECs"
$data[2,1] = "Rspo3"
$data[2,2] = "Lgr5"
$data[2,3] = "MuSCs"
$data[2,4] = 2
$data[2,5] = 0.6666666666666666
$data[2,6] = 0.1043256666666667
$data[2,7] = 0.312977
$data[2,8] = 0.02547563162231953
$data[2,9] = 0.02547563162231953
$data[2,10] = 3
$data[2,11] = 1
$data[2,12] = 3.531223333333333
$data[2,13] = 10.59367
$data[2,14] = 0.5888368093055345
$data[2,15] = 0.5888368093055344
$data[2,16] = 0.3683972283988889
$data[2,17] = 3.31557505559
$data[2,18] = 0.01500098963952981
$data[2,19] = 0.01500098963952981
$data[3,0] = "FAPs"
$data[3,1] = "Rspo3"
$data[3,2] = "Lgr5"
$data[3,3] = "ECs"
$data[3,4] = 3
$data[3,5] = 1
$data[3,6] = 3.990790333333333
$data[3,7] = 11.972371
$data[3,8] = 0.9745243683776804
$data[3,9] = 0.9745243683776804
$data[3,10] = 3
$data[3,11] = 1
$data[3,12] = 0.4967476666666666
$data[3,13] = 1.490243
$data[3,14] = 0.08283342158193596
$data[3,15] = 0.08283342158193596
$data[3,16] = 1.982415786239222
$data[3,17] = 17.841742076153
$data[3,18] = 0.08072318784769826
$data[3,19] = 0.08072318784769826
$data[4,0] = "FAPs"
$data[4,1] = "Rspo3"
$data[4,2] = "Lgr5"
$data[4,3] = "FAPs"
$data[4,4] = 3
$data[4,5] = 1
$data[4,6] = 3.990790333333333
$data[4,7] = 11.972371
$data[4,8] = 0.9745243683776804
$data[4,9] = 0.9745243683776804
$data[4,10] = 3
$data[4,11] = 1
$data[4,12] = 1.968976333333333
$data[4,13] = 5.906929
$data[4,14] = 0.3283297691125296
$data[4,15] = 0.3283297691125296
$data[4,16] = 7.857771717628776
$data[4,17] = 70.719945458659
$data[4,18] = 0.3199653608639775
$data[4,19] = 0.3199653608639775
$data[5,0] = "FAPs"
$data[5,1] = "Rspo3"
$data[5,2] = "Lgr5"
$data[5,3] = "MuSCs"
$data[5,4] = 3
$data[5,5] = 1
$data[5,6] = 3.990790333333333
$data[5,7] = 11.972371
$data[5,8] = 0.9745243683776804
$data[5,9] = 0.9745243683776804
$data[5,10] = 3
$data[5,11] = 1
$data[5,12] = 3.531223333333333
$data[5,13] = 10.59367
$data[5,14] = 0.5888368093055345
$data[5,15] = 0.5888368093055344
$data[5,16] = 14.09237194350778
$data[5,17] = 126.83134749157
$data[5,18] = 0.5738358196660046
$data[5,19] = 0.5738358196660045

$ws.Range("A2:T7").Value2 = $data

Write-Output "done"
